$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add the new "location_type" table column (7th column). This extends the
# table range and the sheet dimension automatically.
$col = $lo.ListColumns.Add()
$ws.Range("G1").Value = "location_type"

# Fill the new column with "Naloxbox" for every existing data row (2-71).
for ($r = 2; $r -le 71; $r++) {
    $ws.Cells.Item($r, 7).Value = "Naloxbox"
}

# Append four new data rows to the table.
$row72 = $lo.ListRows.Add()
$ws.Range("A72").Value = "73 James P Kelly Way, Middletown, NY 10940"
$ws.Range("C72").Value = "Main Office"
$ws.Range("E72").Value = -74.433728900000006
$ws.Range("F72").Value = 41.424500500000001
$ws.Range("G72").Value = "Naloxbox"

$row73 = $lo.ListRows.Add()
$ws.Range("A73").Value = "100 Leprechaun Ln, New Windsor, NY 12553"
$ws.Range("B73").Value = "Hudson House"
$ws.Range("E73").Value = -74.035192699999996
$ws.Range("F73").Value = 41.4921261
$ws.Range("G73").Value = "Naloxbox"

$row74 = $lo.ListRows.Add()
$ws.Range("A74").Value = "123 Pike St, Port Jervis, NY 12771"
$ws.Range("B74").Value = "Hudson House "
$ws.Range("C74").Value = "Suite 208"
$ws.Range("E74").Value = -74.691874299999995
$ws.Range("F74").Value = 41.375512000000001
$ws.Range("G74").Value = "Naloxbox"

$row75 = $lo.ListRows.Add()
$ws.Range("A75").Value = "10 Little Britain Rd, Newburgh, NY 12550"
$ws.Range("E75").Value = -74.032155099999997
$ws.Range("F75").Value = 41.500152
$ws.Range("G75").Value = "Naloxbox"

# Match the saved view state (scroll position / selection) from the diff.
$ws.Range("F78").Select()
